$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G (covid shifts to H)
$ws.Columns.Item(7).Insert()

# New column header
$ws.Range("G1").Value = "off_trade_visibility"

# New off_trade_visibility values for rows 2..158
$values = @(11,10,7,7,10,13,9,9,7,18,36,21,18,8,16,5,5,10,12,18,14,20,28,46,68,57,35,19,14,12,14,10,15,16,13,14,16,15,14,21,13,14,14,9,10,7,11,6,7,9,7,7,6,7,6,6,8,5,4,5,4,7,21,15,9,4,8,5,5,3,3,10,14,19,18,27,23,25,27,14,9,8,6,9,15,12,7,7,9,4,18,16,11,10,10.5,11,10.5,9.5,8,5,6,5,9,8,4,6,6,4,3,3,3,4,4,5,8,7,8,9,6,11,9,7,5,7,10,18,34,27,36,60,54,15,16,17,10,9,18,14,11,9,13,8,17,22,11,9,6,7,13,11,12,12,9,10,9,5,7)
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $values[$i]
}
